$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -532
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -766
$ws.Range("N23").ClearContents()

$ws.Range("H40").Value = 2436.3157
$ws.Range("I40").Value = 2179
$ws.Range("J40").Value = 2722.2222
$ws.Range("K40").Value = 2179
$ws.Range("L40").Value = 2722.2222
$ws.Range("M40").Value = -2004
$ws.Range("N40").Value = -3072.2222

$ws.Range("H106").Value = 38187.562
$ws.Range("I106").Value = 38000.07
$ws.Range("J106").Value = 39500
$ws.Range("K106").Value = 38000.07
$ws.Range("L106").Value = 39500
$ws.Range("M106").Value = -37369.07
$ws.Range("N106").Value = -40762

$ws.Range("H135").Value = 927.9048
$ws.Range("I135").Value = 857.05
$ws.Range("K135").Value = 7713.45
$ws.Range("M135").Value = -5178.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3100.75
$ws.Range("I2").Value = 2022.8
$ws.Range("J2").Value = 4897.3335
$ws.Range("K2").Value = 2022.8
$ws.Range("L2").Value = 4897.3335
$ws.Range("M2").Value = -1909.8
$ws.Range("N2").Value = -5123.3335

$ws.Range("H32").Value = 9150.048000000001
$ws.Range("I32").Value = 6303.2354
$ws.Range("K32").Value = 6303.2354
$ws.Range("M32").Value = -6016.2354

$ws.Range("H45").Value = 1999.2858
$ws.Range("I45").Value = 1999.2858
$ws.Range("K45").Value = 1999.2858
$ws.Range("M45").Value = -1622.2858

$ws.Range("H116").Value = 3100.75
$ws.Range("I116").Value = 2022.8
$ws.Range("J116").Value = 4897.3335
$ws.Range("K116").Value = 2022.8
$ws.Range("L116").Value = 4897.3335
$ws.Range("M116").Value = 271.2
$ws.Range("N116").Value = -9485.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3100.75
$ws.Range("I3").Value = 2022.8
$ws.Range("J3").Value = 4897.3335
$ws.Range("K3").Value = 2022.8
$ws.Range("L3").Value = 4897.3335
$ws.Range("M3").Value = -1908.8
$ws.Range("N3").Value = -5125.3335

$ws.Range("H5").Value = 1646.3334
$ws.Range("I5").Value = 1646.3334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1646.3334
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1533.3334
$ws.Range("N5").ClearContents()

$ws.Range("H134").Value = 2237.077
$ws.Range("I134").Value = 2189.3635
$ws.Range("K134").Value = 6568.0905
$ws.Range("M134").Value = -4033.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2127.389
$ws.Range("I58").Value = 1270
$ws.Range("K58").Value = 1270
$ws.Range("M58").Value = -1067

$ws.Range("H60").Value = 13120.723
$ws.Range("I60").Value = 10951.412
$ws.Range("J60").Value = 49999
$ws.Range("K60").Value = 10951.412
$ws.Range("L60").Value = 49999
$ws.Range("M60").Value = -10440.412
$ws.Range("N60").Value = -51021

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H136").Value = 2127.389
$ws.Range("I136").Value = 1270
$ws.Range("K136").Value = 3810
$ws.Range("M136").Value = -1260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 824.2857
$ws.Range("I14").Value = 824.2857
$ws.Range("K14").Value = 2472.8571
$ws.Range("M14").Value = -2299.8571

$ws.Range("H117").Value = 1656.25
$ws.Range("I117").Value = 409.66666
$ws.Range("J117").Value = 2404.2
$ws.Range("K117").Value = 1228.99998
$ws.Range("L117").Value = 7212.599999999999
$ws.Range("M117").Value = 2213.00002
$ws.Range("N117").Value = -14096.6

$ws.Range("H134").Value = 2705.8
$ws.Range("I134").Value = 2705.8
$ws.Range("K134").Value = 8117.400000000001
$ws.Range("M134").Value = -3047.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 30000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H122").Value = 86463.664
$ws.Range("I122").Value = 3006.6667
$ws.Range("J122").Value = 169920.67
$ws.Range("K122").Value = 9020.000100000001
$ws.Range("L122").Value = 509762.01
$ws.Range("M122").Value = -6570.000100000001
$ws.Range("N122").Value = -514662.01

$ws.Range("H123").Value = 22214.357
$ws.Range("J123").Value = 22214.357
$ws.Range("L123").Value = 22214.357
$ws.Range("N123").Value = -27114.357

$ws.Range("H132").Value = 3149.1
$ws.Range("I132").Value = 2690.0908
$ws.Range("K132").Value = 8070.2724
$ws.Range("M132").Value = -5540.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2394.4119
$ws.Range("I40").Value = 2231.5625
$ws.Range("K40").Value = 2231.5625
$ws.Range("M40").Value = -2095.5625

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H68").Value = 6001.5
$ws.Range("I68").Value = 6000
$ws.Range("K68").Value = 6000
$ws.Range("M68").Value = -5251

$ws.Range("H71").Value = 6001.5
$ws.Range("I71").Value = 6000
$ws.Range("K71").Value = 30000
$ws.Range("M71").Value = -26256

$ws.Range("H88").Value = 27479.25
$ws.Range("J88").Value = 33248.668
$ws.Range("L88").Value = 33248.668
$ws.Range("N88").Value = -34104.668

$ws.Range("H91").Value = 27479.25
$ws.Range("J91").Value = 33248.668
$ws.Range("L91").Value = 33248.668
$ws.Range("N91").Value = -36212.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 15833.333
$ws.Range("J21").Value = 15833.333
$ws.Range("L21").Value = 15833.333
$ws.Range("N21").Value = -16303.333

$ws.Range("H35").Value = 15833.333
$ws.Range("J35").Value = 15833.333
$ws.Range("L35").Value = 15833.333
$ws.Range("N35").Value = -16413.333

$ws.Range("H58").Value = 21450
$ws.Range("I58").Value = 21450
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 21450
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -21142
$ws.Range("N58").ClearContents()

$ws.Range("H62").Value = 7737.5
$ws.Range("J62").Value = 7700
$ws.Range("L62").Value = 7700
$ws.Range("N62").Value = -8948

$ws.Range("H65").Value = 7737.5
$ws.Range("J65").Value = 7700
$ws.Range("L65").Value = 38500
$ws.Range("N65").Value = -44740

$ws.Range("H82").Value = 31000
$ws.Range("J82").Value = 31000
$ws.Range("L82").Value = 31000
$ws.Range("N82").Value = -31766

$ws.Range("H85").Value = 31000
$ws.Range("J85").Value = 31000
$ws.Range("L85").Value = 31000
$ws.Range("N85").Value = -33652

$ws.Range("H136").Value = 1155.4615
$ws.Range("I136").Value = 1125.8334
$ws.Range("J136").Value = 1511
$ws.Range("K136").Value = 3377.5002
$ws.Range("L136").Value = 4533
$ws.Range("M136").Value = -827.5001999999999
$ws.Range("N136").Value = -9633
